$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "C1qa"
$ws.Range("C2").Value = "Cspg4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 141.70809
$ws.Range("H2").Value = 425.12427
$ws.Range("I2").Value = 0.4270657810795758
$ws.Range("J2").Value = 0.4270657810795759
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.162471666666667
$ws.Range("N2").Value = 3.487415
$ws.Range("O2").Value = 0.05948594836865682
$ws.Range("P2").Value = 0.05948594836865682
$ws.Range("Q2").Value = 164.73163956245
$ws.Range("R2").Value = 1482.58475606205
$ws.Range("S2").Value = 0.02540441300331974
$ws.Range("T2").Value = 0.02540441300331975

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "C1qa"
$ws.Range("C3").Value = "Cspg4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 141.70809
$ws.Range("H3").Value = 425.12427
$ws.Range("I3").Value = 0.4270657810795758
$ws.Range("J3").Value = 0.4270657810795759
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.360490666666667
$ws.Range("N3").Value = 16.081472
$ws.Range("O3").Value = 0.2743067897236206
$ws.Range("P3").Value = 0.2743067897236206
$ws.Range("Q3").Value = 759.62489383616
$ws.Range("R3").Value = 6836.624044525441
$ws.Range("S3").Value = 0.117147043408749
$ws.Range("T3").Value = 0.117147043408749

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "C1qa"
$ws.Range("C4").Value = "Cspg4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 141.70809
$ws.Range("H4").Value = 425.12427
$ws.Range("I4").Value = 0.4270657810795758
$ws.Range("J4").Value = 0.4270657810795759
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4960333333333334
$ws.Range("N4").Value = 1.4881
$ws.Range("O4").Value = 0.02538299564789342
$ws.Range("P4").Value = 0.02538299564789342
$ws.Range("Q4").Value = 70.29193624300001
$ws.Range("R4").Value = 632.6274261870001
$ws.Range("S4").Value = 0.01084020886250708
$ws.Range("T4").Value = 0.01084020886250708

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "C1qa"
$ws.Range("C5").Value = "Cspg4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 141.70809
$ws.Range("H5").Value = 425.12427
$ws.Range("I5").Value = 0.4270657810795758
$ws.Range("J5").Value = 0.4270657810795759
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1237266666666667
$ws.Range("N5").Value = 0.37118
$ws.Range("O5").Value = 0.006331335477847643
$ws.Range("P5").Value = 0.006331335477847643
$ws.Range("Q5").Value = 17.5330696154
$ws.Range("R5").Value = 157.7976265386
$ws.Range("S5").Value = 0.002703896731123833
$ws.Range("T5").Value = 0.002703896731123834

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "C1qa"
$ws.Range("C6").Value = "Cspg4"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 141.70809
$ws.Range("H6").Value = 425.12427
$ws.Range("I6").Value = 0.4270657810795758
$ws.Range("J6").Value = 0.4270657810795759
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1821596666666666
$ws.Range("N6").Value = 0.5464789999999999
$ws.Range("O6").Value = 0.009321466352170649
$ws.Range("P6").Value = 0.009321466352170649
$ws.Range("Q6").Value = 25.81349843836999
$ws.Range("R6").Value = 232.32148594533
$ws.Range("S6").Value = 0.003980879308496743
$ws.Range("T6").Value = 0.003980879308496743

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "C1qa"
$ws.Range("C7").Value = "Cspg4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 141.70809
$ws.Range("H7").Value = 425.12427
$ws.Range("I7").Value = 0.4270657810795758
$ws.Range("J7").Value = 0.4270657810795759
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.217072
$ws.Range("N7").Value = 36.65121600000001
$ws.Range("O7").Value = 0.625171464429811
$ws.Range("P7").Value = 0.6251714644298109
$ws.Range("Q7").Value = 1731.25793851248
$ws.Range("R7").Value = 15581.32144661232
$ws.Range("S7").Value = 0.2669893397653795
$ws.Range("T7").Value = 0.2669893397653795

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "C1qa"
$ws.Range("C8").Value = "Cspg4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 171.783722
$ws.Range("H8").Value = 515.3511659999999
$ws.Range("I8").Value = 0.5177047366363254
$ws.Range("J8").Value = 0.5177047366363255
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.162471666666667
$ws.Range("N8").Value = 3.487415
$ws.Range("O8").Value = 0.05948594836865682
$ws.Range("P8").Value = 0.05948594836865682
$ws.Range("Q8").Value = 199.6937096195433
$ws.Range("R8").Value = 1797.24338657589
$ws.Range("S8").Value = 0.03079615723375753
$ws.Range("T8").Value = 0.03079615723375754

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "C1qa"
$ws.Range("C9").Value = "Cspg4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 171.783722
$ws.Range("H9").Value = 515.3511659999999
$ws.Range("I9").Value = 0.5177047366363254
$ws.Range("J9").Value = 0.5177047366363255
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.360490666666667
$ws.Range("N9").Value = 16.081472
$ws.Range("O9").Value = 0.2743067897236206
$ws.Range("P9").Value = 0.2743067897236206
$ws.Range("Q9").Value = 920.8450384662614
$ws.Range("R9").Value = 8287.605346196351
$ws.Range("S9").Value = 0.1420099243314229
$ws.Range("T9").Value = 0.1420099243314229

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "C1qa"
$ws.Range("C10").Value = "Cspg4"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 171.783722
$ws.Range("H10").Value = 515.3511659999999
$ws.Range("I10").Value = 0.5177047366363254
$ws.Range("J10").Value = 0.5177047366363255
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4960333333333334
$ws.Range("N10").Value = 1.4881
$ws.Range("O10").Value = 0.02538299564789342
$ws.Range("P10").Value = 0.02538299564789342
$ws.Range("Q10").Value = 85.21045223606667
$ws.Range("R10").Value = 766.8940701246
$ws.Range("S10").Value = 0.01314089707693366
$ws.Range("T10").Value = 0.01314089707693366

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "C1qa"
$ws.Range("C11").Value = "Cspg4"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 171.783722
$ws.Range("H11").Value = 515.3511659999999
$ws.Range("I11").Value = 0.5177047366363254
$ws.Range("J11").Value = 0.5177047366363255
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1237266666666667
$ws.Range("N11").Value = 0.37118
$ws.Range("O11").Value = 0.006331335477847643
$ws.Range("P11").Value = 0.006331335477847643
$ws.Range("Q11").Value = 21.25422731065333
$ws.Range("R11").Value = 191.28804579588
$ws.Range("S11").Value = 0.003277762366115338
$ws.Range("T11").Value = 0.003277762366115338

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "C1qa"
$ws.Range("C12").Value = "Cspg4"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 171.783722
$ws.Range("H12").Value = 515.3511659999999
$ws.Range("I12").Value = 0.5177047366363254
$ws.Range("J12").Value = 0.5177047366363255
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1821596666666666
$ws.Range("N12").Value = 0.5464789999999999
$ws.Range("O12").Value = 0.009321466352170649
$ws.Range("P12").Value = 0.009321466352170649
$ws.Range("Q12").Value = 31.29206553827932
$ws.Range("R12").Value = 281.6285898445139
$ws.Range("S12").Value = 0.004825767282914875
$ws.Range("T12").Value = 0.004825767282914876

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "C1qa"
$ws.Range("C13").Value = "Cspg4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 171.783722
$ws.Range("H13").Value = 515.3511659999999
$ws.Range("I13").Value = 0.5177047366363254
$ws.Range("J13").Value = 0.5177047366363255
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 12.217072
$ws.Range("N13").Value = 36.65121600000001
$ws.Range("O13").Value = 0.625171464429811
$ws.Range("P13").Value = 0.6251714644298109
$ws.Range("Q13").Value = 2098.694100101984
$ws.Range("R13").Value = 18888.24690091786
$ws.Range("S13").Value = 0.3236542283451812
$ws.Range("T13").Value = 0.3236542283451812

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "C1qa"
$ws.Range("C14").Value = "Cspg4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 18.32613333333333
$ws.Range("H14").Value = 54.9784
$ws.Range("I14").Value = 0.05522948228409861
$ws.Range("J14").Value = 0.05522948228409861
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.162471666666667
$ws.Range("N14").Value = 3.487415
$ws.Range("O14").Value = 0.05948594836865682
$ws.Range("P14").Value = 0.05948594836865682
$ws.Range("Q14").Value = 21.30361075955556
$ws.Range("R14").Value = 191.732496836
$ws.Range("S14").Value = 0.003285378131579536
$ws.Range("T14").Value = 0.003285378131579536

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "C1qa"
$ws.Range("C15").Value = "Cspg4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 18.32613333333333
$ws.Range("H15").Value = 54.9784
$ws.Range("I15").Value = 0.05522948228409861
$ws.Range("J15").Value = 0.05522948228409861
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.360490666666667
$ws.Range("N15").Value = 16.081472
$ws.Range("O15").Value = 0.2743067897236206
$ws.Range("P15").Value = 0.2743067897236206
$ws.Range("Q15").Value = 98.23706668942224
$ws.Range("R15").Value = 884.1336002048001
$ws.Range("S15").Value = 0.01514982198344867
$ws.Range("T15").Value = 0.01514982198344867

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "C1qa"
$ws.Range("C16").Value = "Cspg4"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 18.32613333333333
$ws.Range("H16").Value = 54.9784
$ws.Range("I16").Value = 0.05522948228409861
$ws.Range("J16").Value = 0.05522948228409861
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.4960333333333334
$ws.Range("N16").Value = 1.4881
$ws.Range("O16").Value = 0.02538299564789342
$ws.Range("P16").Value = 0.02538299564789342
$ws.Range("Q16").Value = 9.090373004444446
$ws.Range("R16").Value = 81.81335704000001
$ws.Range("S16").Value = 0.001401889708452682
$ws.Range("T16").Value = 0.001401889708452682

# Row 17
$ws.Range("A17").Value = "Neutro"
$ws.Range("B17").Value = "C1qa"
$ws.Range("C17").Value = "Cspg4"
$ws.Range("D17").Value = "M2"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 18.32613333333333
$ws.Range("H17").Value = 54.9784
$ws.Range("I17").Value = 0.05522948228409861
$ws.Range("J17").Value = 0.05522948228409861
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1237266666666667
$ws.Range("N17").Value = 0.37118
$ws.Range("O17").Value = 0.006331335477847643
$ws.Range("P17").Value = 0.006331335477847643
$ws.Range("Q17").Value = 2.267431390222222
$ws.Range("R17").Value = 20.406882512
$ws.Range("S17").Value = 0.0003496763806084714
$ws.Range("T17").Value = 0.0003496763806084714

# Row 18
$ws.Range("A18").Value = "Neutro"
$ws.Range("B18").Value = "C1qa"
$ws.Range("C18").Value = "Cspg4"
$ws.Range("D18").Value = "Neutro"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 18.32613333333333
$ws.Range("H18").Value = 54.9784
$ws.Range("I18").Value = 0.05522948228409861
$ws.Range("J18").Value = 0.05522948228409861
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.1821596666666666
$ws.Range("N18").Value = 0.5464789999999999
$ws.Range("O18").Value = 0.009321466352170649
$ws.Range("P18").Value = 0.009321466352170649
$ws.Range("Q18").Value = 3.338282339288889
$ws.Range("R18").Value = 30.0445410536
$ws.Range("S18").Value = 0.0005148197607590301
$ws.Range("T18").Value = 0.0005148197607590301

# Row 19
$ws.Range("A19").Value = "Neutro"
$ws.Range("B19").Value = "C1qa"
$ws.Range("C19").Value = "Cspg4"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 18.32613333333333
$ws.Range("H19").Value = 54.9784
$ws.Range("I19").Value = 0.05522948228409861
$ws.Range("J19").Value = 0.05522948228409861
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 12.217072
$ws.Range("N19").Value = 36.65121600000001
$ws.Range("O19").Value = 0.625171464429811
$ws.Range("P19").Value = 0.6251714644298109
$ws.Range("Q19").Value = 223.8916904149334
$ws.Range("R19").Value = 2015.0252137344
$ws.Range("S19").Value = 0.03452789631925023
$ws.Range("T19").Value = 0.03452789631925023

